$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column cells that change below are re-written as text so that
# values like "324.40" or "28.884.81" keep their original formatting
# instead of being reinterpreted as numbers (which would drop trailing
# zeros / mis-parse the "thousand dot" notation used in this sheet).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.884.81"
$ws.Range("E2").Value = "  -0.24%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.917.45"
$ws.Range("E3").Value = "  +0.73%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.40"
$ws.Range("E5").Value = "  +0.11%  "

# Row 6
$ws.Range("E6").Value = "  +0.07%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4566"
$ws.Range("E7").Value = "  -0.60%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3802"
$ws.Range("E8").Value = "  -0.25%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07751"
$ws.Range("E9").Value = "  +0.68%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9772"
$ws.Range("E10").Value = "  -0.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.27"
$ws.Range("E11").Value = "  +1.30%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.910.36"
$ws.Range("E12").Value = "  +0.91%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.696"
$ws.Range("E13").Value = "  +0.46%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.968"
$ws.Range("E14").Value = "  +0.16%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06984"
$ws.Range("E15").Value = "  -1.15%  "

# Row 16
$ws.Range("E16").Value = "  -0.06%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.38"
$ws.Range("E17").Value = "  +0.74%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009489"
$ws.Range("E18").Value = "  -0.38%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.65"
$ws.Range("E19").Value = "  -0.33%  "

# Row 20
$ws.Range("E20").Value = "  -0.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.889.27"
$ws.Range("E21").Value = "  -0.17%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.339"
$ws.Range("E22").Value = "  +0.30%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.07"
$ws.Range("E23").Value = "  +1.66%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.143.20"
$ws.Range("E24").Value = "  +0.60%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.060"
$ws.Range("E25").Value = "  -1.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.69"
$ws.Range("E26").Value = "  +0.40%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.98"
$ws.Range("E27").Value = "  -0.59%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.630"
$ws.Range("E28").Value = "  +0.66%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.58"
$ws.Range("E29").Value = "  +0.08%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.850"
$ws.Range("E30").Value = "  +0.10%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09287"
$ws.Range("E31").Value = "  +0.02%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8705"
$ws.Range("E32").Value = "  +1.25%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.101"
$ws.Range("E33").Value = "  +0.52%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.243"
$ws.Range("E34").Value = "  -0.40%  "

# Row 35
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05703"
$ws.Range("E36").Value = "  +0.27%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.149"
$ws.Range("E37").Value = "  +0.27%  "

# Row 38
$ws.Range("E38").Value = "  +0.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02041"
$ws.Range("E39").Value = "  +0.49%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.059"
$ws.Range("E40").Value = "  +11.51%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.488"
$ws.Range("E41").Value = "  +0.44%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5511"
$ws.Range("E42").Value = "  +0.07%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1756"
$ws.Range("E43").Value = "  +0.06%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.341"
$ws.Range("E44").Value = "  +0.98%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002840"
$ws.Range("E45").Value = "  +14.99%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.165"
$ws.Range("E46").Value = "  +3.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5161"
$ws.Range("E47").Value = "  -0.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06924"
$ws.Range("E48").Value = "  +1.50%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.12"
$ws.Range("E49").Value = "  -1.82%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.82"
$ws.Range("E50").Value = "  -0.43%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.768"
$ws.Range("E51").Value = "  -0.41%  "
